# Apply the psMCU Encodings.xlsx edits:
#  - "Inst Decoding" sheet: rename several mnemonics (Save/Load rework) and
#    fix "KOMP" -> "COMP" typo.
#  - "Fixed Peripheral Registers" sheet: replace the placeholder 0/0/0 row
#    with the real flag labels A=B / B=0 / A=0.

$wb = $excel.ActiveWorkbook

$wsInst = $wb.Worksheets.Item("Inst Decoding")
$wsInst.Range("B7").Value = "SV[A/B]"
$wsInst.Range("B8").Value = "LD[A/B]"
$wsInst.Range("B9").Value = "LIT"
$wsInst.Range("B12").Value = "SVD[P/M]"
$wsInst.Range("B13").Value = "LDD[P/M]"

$wsRegs = $wb.Worksheets.Item("Fixed Peripheral Registers")
$wsRegs.Range("C3").Value = "A=0"
$wsRegs.Range("B3").Value = "B=0"
$wsRegs.Range("A3").Value = "A=B"

$wsInst.Range("A16").Value = "COMP B -> B"
$wsInst.Range("B16").Value = "COMPB"
